$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "28.306.80"
$ws.Range("E2").Value = "  +0.39%  "

$ws.Range("D3").Value = "1.872.03"
$ws.Range("E3").Value = "  +3.67%  "

$ws.Range("D4").Value = "1.002"
$ws.Range("E4").Value = "  +0.10%  "

$ws.Range("D5").Value = "311.96"
$ws.Range("E5").Value = "  -0.22%  "

$ws.Range("E6").Value = "  +0.10%  "

$ws.Range("D7").Value = "0.5020"
$ws.Range("E7").Value = "  -2.28%  "

$ws.Range("D8").Value = "0.3961"
$ws.Range("E8").Value = "  -0.18%  "

$ws.Range("D9").Value = "0.09907"
$ws.Range("E9").Value = "  +27.02%  "

$ws.Range("D10").Value = "1.124"
$ws.Range("E10").Value = "  +1.55%  "

$ws.Range("D11").Value = "41.33"
$ws.Range("E11").Value = "  +1.12%  "

$ws.Range("D12").Value = "6.492"
$ws.Range("E12").Value = "  +2.18%  "

$ws.Range("D13").Value = "21.02"
$ws.Range("E13").Value = "  +2.99%  "

$ws.Range("D14").Value = "1.870.26"
$ws.Range("E14").Value = "  +3.81%  "

$ws.Range("E15").Value = "  +0.12%  "

$ws.Range("D16").Value = "7.386"
$ws.Range("E16").Value = "  +1.14%  "

$ws.Range("D17").Value = "0.00001143"
$ws.Range("E17").Value = "  +6.01%  "

$ws.Range("D18").Value = "93.92"
$ws.Range("E18").Value = "  +1.00%  "

$ws.Range("D19").Value = "0.06693"
$ws.Range("E19").Value = "  +1.89%  "

$ws.Range("E20").Value = "  +0.16%  "

$ws.Range("E21").Value = "  +0.73%  "

$ws.Range("D22").Value = "6.109"
$ws.Range("E22").Value = "  +1.86%  "

$ws.Range("D23").Value = "28.382.12"
$ws.Range("E23").Value = "  +0.55%  "

$ws.Range("D24").Value = "11.35"
$ws.Range("E24").Value = "  +2.01%  "

$ws.Range("D25").Value = "2.262"
$ws.Range("E25").Value = "  +2.23%  "

$ws.Range("D26").Value = "2.523"
$ws.Range("E26").Value = "  +3.67%  "

$ws.Range("E27").Value = "  +3.62%  "

$ws.Range("D28").Value = "2.084.21"
$ws.Range("E28").Value = "  +3.48%  "

$ws.Range("D29").Value = "157.92"
$ws.Range("E29").Value = "  -1.62%  "

$ws.Range("D30").Value = "127.77"
$ws.Range("E30").Value = "  -0.08%  "

$ws.Range("E31").Value = "  -2.95%  "

$ws.Range("E32").Value = "  +0.43%  "

$ws.Range("D33").Value = "5.650"
$ws.Range("E33").Value = "  +1.56%  "

$ws.Range("D34").Value = "3.611"
$ws.Range("E34").Value = "  -1.17%  "

$ws.Range("D35").Value = "0.06812"
$ws.Range("E35").Value = "  -4.95%  "

$ws.Range("D36").Value = "9.365"
$ws.Range("E36").Value = "  +2.97%  "

$ws.Range("D37").Value = "0.02396"
$ws.Range("E37").Value = "  +2.07%  "

$ws.Range("D38").Value = "0.2194"
$ws.Range("E38").Value = "  +0.98%  "

$ws.Range("D39").Value = "5.026"
$ws.Range("E39").Value = "  -0.15%  "

$ws.Range("D40").Value = "11.50"
$ws.Range("E40").Value = "  -0.24%  "

$ws.Range("D41").Value = "0.6293"
$ws.Range("E41").Value = "  +2.15%  "

$ws.Range("D42").Value = "1.176"
$ws.Range("E42").Value = "  +2.02%  "

$ws.Range("E43").Value = "  +0.08%  "

$ws.Range("D44").Value = "13.43"
$ws.Range("E44").Value = "  +1.89%  "

$ws.Range("D45").Value = "0.6008"
$ws.Range("E45").Value = "  +0.92%  "

$ws.Range("B46").Value = "WEMIXTOKEN"
$ws.Range("C46").Value = "https://coinranking.com/coin/08CsQa-Ov+wemixtoken-wemix"
$ws.Range("D46").Value = "1.283"
$ws.Range("E46").Value = "  -2.03%  "

$ws.Range("B47").Value = "PancakeSwap"
$ws.Range("C47").Value = "https://coinranking.com/coin/ncYFcP709+pancakeswap-cake"
$ws.Range("D47").Value = "3.677"
$ws.Range("E47").Value = "  -1.63%  "

$ws.Range("D48").Value = "125.16"
$ws.Range("E48").Value = "  +0.34%  "

$ws.Range("D49").Value = "1.993"
$ws.Range("E49").Value = "  +4.12%  "

$ws.Range("D50").Value = "1.200"
$ws.Range("E50").Value = "  -1.14%  "

$ws.Range("B51").Value = "Cronos"
$ws.Range("C51").Value = "https://coinranking.com/coin/65PHZTpmE55b+cronos-cro"
$ws.Range("D51").Value = "0.06852"
$ws.Range("E51").Value = "  +0.81%  "
